$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New match rows appended below the existing data (rows 13-23).
# Columns: A venue, B date, C result, D ownTeam, E oppTeam, F batsman,
#          G totalRuns, H totalBalls, I total4s, J total6s, K sr
$newRows = @(
    @(" Dubai (DSC)", " October 14 2020", "Capitals won by 13 runs", "Rajasthan Royals", "Delhi Capitals", "Rahul Tewatia ", "14", "18", "1", "0", "77.77"),
    @(" Dubai (DSC)", " October 22 2020", "Sunrisers won by 8 wickets (with 11 balls remaining)", "Rajasthan Royals", "Sunrisers Hyderabad", "Rahul Tewatia ", "2", "3", "0", "0", "66.66"),
    @(" Sharjah", " September 22 2020", "Royals won by 16 runs", "Rajasthan Royals", "Chennai Super Kings", "Rahul Tewatia ", "10", "8", "1", "0", "125.00"),
    @(" Dubai (DSC)", " September 30 2020", "KKR won by 37 runs", "Rajasthan Royals", "Kolkata Knight Riders", "Rahul Tewatia ", "14", "10", "0", "1", "140.00"),
    @(" Dubai (DSC)", " November 01 2020", "KKR won by 60 runs", "Rajasthan Royals", "Kolkata Knight Riders", "Rahul Tewatia ", "31", "27", "2", "1", "114.81"),
    @(" Sharjah", " September 27 2020", "Royals won by 4 wickets (with 3 balls remaining)", "Rajasthan Royals", "Kings XI Punjab", "Rahul Tewatia ", "53", "31", "0", "7", "170.96"),
    @(" Dubai (DSC)", " October 17 2020", "RCB won by 7 wickets (with 2 balls remaining)", "Rajasthan Royals", "Royal Challengers Bangalore", "Rahul Tewatia ", "19", "11", "1", "1", "172.72"),
    @(" Dubai (DSC)", " October 11 2020", "Royals won by 5 wickets (with 1 ball remaining)", "Rajasthan Royals", "Sunrisers Hyderabad", "Rahul Tewatia ", "45", "28", "4", "2", "160.71"),
    @(" Abu Dhabi", " October 06 2020", "Mumbai won by 57 runs", "Rajasthan Royals", "Mumbai Indians", "Rahul Tewatia ", "5", "6", "1", "0", "83.33"),
    @(" Abu Dhabi", " October 03 2020", "RCB won by 8 wickets (with 5 balls remaining)", "Rajasthan Royals", "Royal Challengers Bangalore", "Rahul Tewatia ", "24", "12", "0", "3", "200.00"),
    @(" Sharjah", " October 09 2020", "Capitals won by 46 runs", "Rajasthan Royals", "Delhi Capitals", "Rahul Tewatia ", "38", "29", "3", "2", "131.03")
)

$startRow = 13
$numericCols = @(7, 8, 9, 10, 11)

for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $rowValues = $newRows[$i]
    for ($c = 1; $c -le $rowValues.Count; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        if ($numericCols -contains $c) {
            # Keep these numeric-looking values stored as text, matching the source data.
            $cell.NumberFormat = "@"
        }
        $cell.Value = $rowValues[$c - 1]
    }
}
